$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new H1 cell,
# then set its text, so the new column reuses the existing header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Save values for rows 2-8 (1 indicates a save, 0 otherwise)
$saveValues = @(0, 0, 0, 0, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
